$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 163.89473
$ws.Range("J17").Value = 163.89473
$ws.Range("L17").Value = 491.6841900000001
$ws.Range("N17").Value = -827.6841900000001

$ws.Range("H19").Value = 1580.1666
$ws.Range("I19").Value = 800
$ws.Range("J19").Value = 1970.25
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 1970.25
$ws.Range("M19").Value = -625
$ws.Range("N19").Value = -2320.25

$ws.Range("H100").Value = 15153270
$ws.Range("I100").Value = 18519580
$ws.Range("J100").Value = 4875
$ws.Range("K100").Value = 18519580
$ws.Range("L100").Value = 4875
$ws.Range("M100").Value = -18519039
$ws.Range("N100").Value = -5957

$ws.Range("H125").Value = 779.9091
$ws.Range("I125").Value = 290
$ws.Range("J125").Value = 888.7778
$ws.Range("K125").Value = 2610
$ws.Range("L125").Value = 7999.000199999999
$ws.Range("M125").Value = -150
$ws.Range("N125").Value = -12919.0002

$ws.Range("H131").Value = 2691.0557
$ws.Range("I131").Value = 831.7273
$ws.Range("J131").Value = 5612.857
$ws.Range("K131").Value = 2495.1819
$ws.Range("L131").Value = 16838.571
$ws.Range("M131").Value = 2544.8181
$ws.Range("N131").Value = -26918.571

$ws.Range("H135").Value = 2106.1177
$ws.Range("I135").Value = 2438.4285
$ws.Range("J135").Value = 555.3333
$ws.Range("K135").Value = 21945.8565
$ws.Range("L135").Value = 4997.9997
$ws.Range("M135").Value = -19410.8565
$ws.Range("N135").Value = -10067.9997

$ws.Range("H137").Value = 1040.225
$ws.Range("J137").Value = 1249
$ws.Range("L137").Value = 3747
$ws.Range("N137").Value = -8847

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 305943.3
$ws.Range("I32").Value = 2642.8735
$ws.Range("J32").Value = 3604335.5
$ws.Range("K32").Value = 2642.8735
$ws.Range("L32").Value = 3604335.5
$ws.Range("M32").Value = -2355.8735
$ws.Range("N32").Value = -3604909.5

$ws.Range("H45").Value = 2656.25
$ws.Range("I45").Value = 2614.2856
$ws.Range("J45").Value = 2736.3635
$ws.Range("K45").Value = 2614.2856
$ws.Range("L45").Value = 2736.3635
$ws.Range("M45").Value = -2237.2856
$ws.Range("N45").Value = -3490.3635

$ws.Range("H61").Value = 1838.3726
$ws.Range("I61").Value = 1615.3846
$ws.Range("J61").Value = 2563.0833
$ws.Range("K61").Value = 1615.3846
$ws.Range("L61").Value = 2563.0833
$ws.Range("M61").Value = -1403.3846
$ws.Range("N61").Value = -2987.0833

$ws.Range("H74").Value = 669.5135
$ws.Range("I74").Value = 552.875
$ws.Range("J74").Value = 1416
$ws.Range("K74").Value = 552.875
$ws.Range("L74").Value = 1416
$ws.Range("M74").Value = 321.125
$ws.Range("N74").Value = -3164

$ws.Range("H77").Value = 669.5135
$ws.Range("I77").Value = 552.875
$ws.Range("J77").Value = 1416
$ws.Range("K77").Value = 2764.375
$ws.Range("L77").Value = 7080
$ws.Range("M77").Value = 1603.625
$ws.Range("N77").Value = -15816

$ws.Range("H136").Value = 1838.3726
$ws.Range("I136").Value = 1615.3846
$ws.Range("J136").Value = 2563.0833
$ws.Range("K136").Value = 4846.1538
$ws.Range("L136").Value = 7689.249899999999
$ws.Range("M136").Value = -2296.1538
$ws.Range("N136").Value = -12789.2499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 19241.111
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 19241.111
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 19241.111
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -19767.111

$ws.Range("H86").Value = 65087.312
$ws.Range("I86").Value = 2824.3333
$ws.Range("J86").Value = 251876.25
$ws.Range("K86").Value = 2824.3333
$ws.Range("L86").Value = 251876.25
$ws.Range("M86").Value = -1701.3333
$ws.Range("N86").Value = -254122.25

$ws.Range("H89").Value = 65087.312
$ws.Range("I89").Value = 2824.3333
$ws.Range("J89").Value = 251876.25
$ws.Range("K89").Value = 14121.6665
$ws.Range("L89").Value = 1259381.25
$ws.Range("M89").Value = -8505.666499999999
$ws.Range("N89").Value = -1270613.25

$ws.Range("H94").Value = 1829
$ws.Range("I94").Value = 1541
$ws.Range("J94").Value = 2405
$ws.Range("K94").Value = 1541
$ws.Range("L94").Value = 2405
$ws.Range("M94").Value = -1090
$ws.Range("N94").Value = -3307

$ws.Range("H105").Value = 2962.3333
$ws.Range("I105").Value = 2575.5
$ws.Range("K105").Value = 2575.5
$ws.Range("M105").Value = -828.5

$ws.Range("H121").Value = 19241.111
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 19241.111
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 19241.111
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -22735.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1112521.6
$ws.Range("I31").Value = 1906222.1
$ws.Range("J31").Value = 1340.84
$ws.Range("K31").Value = 1906222.1
$ws.Range("L31").Value = 1340.84
$ws.Range("M31").Value = -1905927.1
$ws.Range("N31").Value = -1930.84

$ws.Range("H34").Value = 1112521.6
$ws.Range("I34").Value = 1906222.1
$ws.Range("J34").Value = 1340.84
$ws.Range("K34").Value = 1906222.1
$ws.Range("L34").Value = 1340.84
$ws.Range("M34").Value = -1906020.1
$ws.Range("N34").Value = -1744.84

$ws.Range("H53").Value = 21047
$ws.Range("J53").Value = 21047
$ws.Range("L53").Value = 21047
$ws.Range("N53").Value = -22261

$ws.Range("H100").Value = 38500
$ws.Range("J100").Value = 38500
$ws.Range("L100").Value = 38500
$ws.Range("N100").Value = -40664

$ws.Range("H122").Value = 7143639.5
$ws.Range("I122").Value = 653
$ws.Range("K122").Value = 1959
$ws.Range("M122").Value = 491

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 160
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 300
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = -160
$ws.Range("N11").Value = -880

$ws.Range("H16").Value = 2150
$ws.Range("J16").Value = 4000
$ws.Range("L16").Value = 12000
$ws.Range("N16").Value = -12346

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H68").Value = 1173
$ws.Range("I68").Value = 707.06665
$ws.Range("J68").Value = 1476.8695
$ws.Range("K68").Value = 2121.19995
$ws.Range("L68").Value = 4430.6085
$ws.Range("M68").Value = -1310.19995
$ws.Range("N68").Value = -6052.6085

$ws.Range("H71").Value = 1173
$ws.Range("I71").Value = 707.06665
$ws.Range("J71").Value = 1476.8695
$ws.Range("K71").Value = 6363.59985
$ws.Range("L71").Value = 13291.8255
$ws.Range("M71").Value = -2307.59985
$ws.Range("N71").Value = -21403.8255

$ws.Range("H107").Value = 765.61365
$ws.Range("I107").Value = 538.46155
$ws.Range("J107").Value = 1093.7222
$ws.Range("K107").Value = 1615.38465
$ws.Range("L107").Value = 3281.1666
$ws.Range("M107").Value = 304.61535
$ws.Range("N107").Value = -7121.1666

$ws.Range("H131").Value = 5377357.5
$ws.Range("I131").Value = 818.3
$ws.Range("J131").Value = 6025133.5
$ws.Range("K131").Value = 2454.9
$ws.Range("L131").Value = 18075400.5
$ws.Range("M131").Value = 2585.1
$ws.Range("N131").Value = -18085480.5

$ws.Range("H138").Value = 1252.6
$ws.Range("I138").Value = 557.5
$ws.Range("J138").Value = 4033
$ws.Range("K138").Value = 1672.5
$ws.Range("L138").Value = 12099
$ws.Range("M138").Value = 3467.5
$ws.Range("N138").Value = -22379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14429536
$ws.Range("I70").Value = 21640166
$ws.Range("K70").Value = 21640166
$ws.Range("M70").Value = -21639896

$ws.Range("H73").Value = 14429536
$ws.Range("I73").Value = 21640166
$ws.Range("K73").Value = 21640166
$ws.Range("M73").Value = -21639230

$ws.Range("H95").Value = 19172
$ws.Range("J95").Value = 19172
$ws.Range("L95").Value = 19172
$ws.Range("N95").Value = -24664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1775.25
$ws.Range("I16").Value = 1967
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 1967
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -1797
$ws.Range("N16").Value = -1540

$ws.Range("H22").Value = 535.1905
$ws.Range("I22").Value = 462.6154
$ws.Range("J22").Value = 653.125
$ws.Range("K22").Value = 462.6154
$ws.Range("L22").Value = 653.125
$ws.Range("M22").Value = -167.6154
$ws.Range("N22").Value = -1243.125

$ws.Range("H27").Value = 535.1905
$ws.Range("I27").Value = 462.6154
$ws.Range("J27").Value = 653.125
$ws.Range("K27").Value = 462.6154
$ws.Range("L27").Value = 653.125
$ws.Range("M27").Value = -355.6154
$ws.Range("N27").Value = -867.125

$ws.Range("H46").Value = 10072.818
$ws.Range("I46").Value = 871.5714
$ws.Range("J46").Value = 26175
$ws.Range("K46").Value = 871.5714
$ws.Range("L46").Value = 26175
$ws.Range("M46").Value = -683.5714
$ws.Range("N46").Value = -26551

$ws.Range("H132").Value = 2166732.8
$ws.Range("I132").Value = 3572699.2
$ws.Range("J132").Value = 3707.1538
$ws.Range("K132").Value = 10718097.6
$ws.Range("L132").Value = 11121.4614
$ws.Range("M132").Value = -10715567.6
$ws.Range("N132").Value = -16181.4614

$ws.Range("H136").Value = 50053276
$ws.Range("I136").Value = 84166.336
$ws.Range("J136").Value = 125006936
$ws.Range("K136").Value = 252499.008
$ws.Range("L136").Value = 375020808
$ws.Range("M136").Value = -249949.008
$ws.Range("N136").Value = -375025908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 620.5
$ws.Range("I113").Value = 620.5
$ws.Range("K113").Value = 1861.5
$ws.Range("M113").Value = 308.5

$ws.Range("H132").Value = 66182140
$ws.Range("I132").Value = 132355430
$ws.Range("J132").Value = 8852.883
$ws.Range("K132").Value = 397066290
$ws.Range("L132").Value = 26558.649
$ws.Range("M132").Value = -397063760
$ws.Range("N132").Value = -31618.649
